$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.690.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.295.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.500'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.494'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.75'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.90%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '16.73'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.03%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.655.29'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.297.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("E18").Value = '  -0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.641.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.56%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0898'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.46'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.38%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.97'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0694'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0996'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("E41").Value = '  -3.63%  '
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.35'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.960.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0280'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.522.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.11%  '
$ws.Range("E51").Value = '  -4.13%  '
